# Cálculo duraciones PE.xlsx - apply the recorded edit
#
# Summary of the change (from the "Diseño" section of the sheet):
#   - Rows 39-42 (raw data rows under "Diseño") get new values.
#   - A brand-new data row is inserted right after the old row 42, containing (12, 12).
#   - The SUM row for the "Diseño" section (previously row 43) moves down to row 44
#     and its formula range grows to include the new row (A36:A43 / B36:B43).
#   - The "Implementación" section header (previously row 44) moves down to row 45.
#   - The "Implementación" data rows (previously 45-52) move down to 46-53, and
#     several of the data values themselves are corrected.
#   - The "Implementación" SUM row (previously row 53) moves down to row 54... but the
#     remaining rows below (Pruebas section, grand total) stay where they were, so the
#     SUM row keeps row number 53 and its range is now B46:B52 (it no longer includes
#     the header row).
#   - The grand-total row (64) formulas are repointed from A43/B43 to A44/B44 (since
#     the "Diseño" subtotal cell moved there).
#
# Because rows below the "Implementación" subtotal are NOT shifted (the "Pruebas"
# header stays at row 54, etc.), this is implemented as direct cell writes rather
# than a literal row-insert/shift, which would have moved everything below too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- "Diseño" section raw data (rows 39-43) ----
$ws.Range("A39").Value = 2
$ws.Range("B39").Value = 5

$ws.Range("A40").Value = 106
$ws.Range("B40").Value = 50

$ws.Range("A41").Value = 2
$ws.Range("B41").Value = 2

$ws.Range("A42").Value = 2
$ws.Range("B42").Value = 2

# New row of data that did not exist before. It reuses the row 43 slot, which used to
# hold the section's SUM formula (bold style) - reset it to a plain data-row look by
# copying the formatting from a neighboring plain data row before writing the values.
$ws.Range("A42:B42").Copy()
$ws.Range("A43:B43").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows("43").AutoFit()

$ws.Range("A43").Value = 12
$ws.Range("B43").Value = 12

# ---- "Diseño" subtotal, now on row 44 ----
# Give it the same look (bold font style) as the other subtotal rows, e.g. row 34,
# and make sure the row reverts to the standard (non-custom) row height.
$ws.Range("A34:B34").Copy()
$ws.Range("A44:B44").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows("44").AutoFit()

$ws.Range("A44").Formula = "=SUM(A36:A43)"
$ws.Range("B44").Formula = "=SUM(B36:B43)"

# ---- "Implementación" header, now on row 45 ----
# Copy the formatting (font/style + row height) from another section header, e.g. row 35.
$ws.Range("A35:B35").Copy()
$ws.Range("A45:B45").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows("45").RowHeight = 26

$ws.Range("A45").Value = "Implementación"
$ws.Range("B45").Value = "Implementación"

# ---- "Implementación" raw data, now on rows 46-52 ----
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = 5

$ws.Range("A47").Value = 5
$ws.Range("B47").Value = 5

$ws.Range("A48").Value = 0
$ws.Range("B48").Value = 5

$ws.Range("A49").Value = 5
$ws.Range("B49").Value = 6

$ws.Range("A50").Value = 0
$ws.Range("B50").Value = 2

$ws.Range("A51").Value = 0
$ws.Range("B51").Value = 2

$ws.Range("A52").Value = 2
$ws.Range("B52").Value = 28

# ---- "Implementación" subtotal, stays on row 53 but formula range shrinks by one
#      row (it no longer spans the header row) and shifts down into the data block ----
$ws.Range("A53").Formula = "=SUM(A46:A52)"
$ws.Range("B53").Formula = "=SUM(B46:B52)"

# ---- Grand total row: repoint references from the old Diseño-subtotal cell (A43/B43)
#      to the new one (A44/B44) ----
$ws.Range("A64").Formula = "=SUM(A18,A34,A44,A53,A61)"
$ws.Range("B64").Formula = "=SUM(B18,B34,B44,B53,B61)"

# ---- Selection / scroll position recorded in the saved view ----
$ws.Range("D52").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
